# Adobe AAM sheet: add a "Segment Lifetime" column and two new test segment rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Adobe AAM")
$ws.Activate()

# Insert a new column before the current "Trait Folder Path" column (E),
# then clone formatting (incl. column width) from what is now column F
# (the old column E) onto the freshly inserted, still-blank column E.
$ws.Columns("E").Insert()
$ws.Range("F1:F2").Copy($ws.Range("E1:E2"))
$ws.Columns("E").ColumnWidth = $ws.Columns("F").ColumnWidth

# New header + "required" note for the inserted column.
$ws.Range("E1").Value = "Segment Lifetime"
$ws.Range("E2").Value = "Add: Required"

# Refresh the trait-folder-path sample values (now in column F) and the
# data-source-name samples (now in column H) for the new test date.
$ws.Range("F3").Value = "/All Traits/TEST20181030/TEST"
$ws.Range("F4").Value = "/All Traits/TEST20181030"
$ws.Range("H3").Value = "test20181030"
$ws.Range("H4").Value = "test20181030"

# Two brand-new example segment rows.
$ws.Range("B3").Value = "Test Segment 1"
$ws.Range("C3").Value = "Test Description 1"
$ws.Range("E3").Value = 90

$ws.Range("B4").Value = "Test Segment 2"
$ws.Range("C4").Value = "Test Description 2"
$ws.Range("E4").Value = 90
$ws.Range("J4").Value = 1

$ws.Range("F11").Select()
